$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46061 -> 46062) for every data row (rows 2 through 96).
for ($r = 2; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
